# Week 4 Plan - mark additional linked-list questions as completed (YES)
# and align column B centered, matching Excel's behaviour when the
# whole column is selected and "Center" alignment is applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center-align the whole of column B (this updates every already-present
# cell in column B that isn't already center-aligned, same as clicking the
# column header and pressing the Center-align button in Excel).
$ws.Columns("B").HorizontalAlignment = -4108

# Mark the newly completed tasks with "YES" in column B.
$ws.Range("B3").Value = "YES"
$ws.Range("B4").Value = "YES"
$ws.Range("B5").Value = "YES"
$ws.Range("B6").Value = "YES"
$ws.Range("B7").Value = "YES"
$ws.Range("B9").Value = "YES"
$ws.Range("B11").Value = "YES"
$ws.Range("B12").Value = "YES"
$ws.Range("B14").Value = "YES"

# Leave the current selection on A10, matching where the author ended up.
$ws.Range("A10").Select() | Out-Null
